$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '39.587.36'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +1.85%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.157.09'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +1.55%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '226.69'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.68%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.621'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.66%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '62.77'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.65%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0841'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.52%  '
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '15.82'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.68%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.473.97'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.67%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '21.71'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.94%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.804'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.60%  '
$ws.Range('E16').Value = '  -0.83%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.162.32'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.56%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '39.528.85'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.58%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '71.55'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.40%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.03'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.86%  '
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '227.49'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  +0.96%  '
$ws.Range('E25').Value = '  -3.75%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '170.82'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.41'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.74%  '
$ws.Range('E28').Value = '  +1.38%  '
$ws.Range('E29').Value = '  +0.60%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '19.61'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('E32').Value = '  +0.41%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.57'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.80%  '
$ws.Range('E34').Value = '  -1.84%  '
$ws.Range('E35').Value = '  -3.32%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0615'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('E37').Value = '  +7.40%  '
$ws.Range('E38').Value = '  -0.40%  '
$ws.Range('B39').Value = 'FTXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.05'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +21.58%  '
$ws.Range('B40').Value = 'BinanceUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.08%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '102.55'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '17.66'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.91%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.513.19'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.12%  '
$ws.Range('E45').Value = '  -0.38%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '7.85'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0920'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.36%  '
$ws.Range('B48').Value = 'HuobiToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.80'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.06%  '
$ws.Range('E49').Value = '  -0.56%  '
$ws.Range('B50').Value = 'TerraClassic'
$ws.Range('C50').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.000190'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +27.71%  '
$ws.Range('B51').Value = 'MXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.98'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.31%  '
